# "Highlighter and report added" - refresh the per-language and per-level
# course-count report figures (Language and Level sheets) to their newly
# computed values. All of these count cells are stored as TEXT in the
# workbook (Apache POI wrote them as shared strings, not numbers), so we
# force a Text number format before writing each value to keep them as
# text rather than letting Excel auto-coerce them to numeric cells.

$wb = $excel.ActiveWorkbook

# ---- "Language" sheet: INDIVIDUAL COURSE COUNT column (B) updates ----
$wsLang = $wb.Worksheets.Item("Language")

$langCounts = [ordered]@{
    "B2"  = "1239"   # English
    "B3"  = "624"    # Russian
    "B4"  = "614"    # Spanish
    "B5"  = "577"    # French
    "B6"  = "536"    # Portuguese (Portugal)
    "B7"  = "480"    # Arabic
    "B8"  = "473"    # Vietnamese
    "B9"  = "472"    # German
    "B10" = "455"    # Italian
    "B11" = "136"    # Chinese (China)
    "B12" = "126"    # Korean
    "B13" = "69"     # Portuguese (Brazil)
    "B24" = "11"     # (new) Portuguese row count
    "B52" = "11681"  # TOTAL COURSES
}

foreach ($addr in $langCounts.Keys) {
    $cell = $wsLang.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $langCounts[$addr]
}

# Row 23/24 language names swap: Urdu now sorts above Portuguese.
$wsLang.Range("A23").Value = "Urdu"
$wsLang.Range("A24").Value = "Portuguese"

# ---- "Level" sheet: INDIVIDUAL LEVEL COURSES column (B) updates ----
$wsLevel = $wb.Worksheets.Item("Level")

$levelCounts = [ordered]@{
    "B2" = "594"   # Intermediate
    "B3" = "584"   # Beginner
    "B4" = "155"   # Mixed
    "B5" = "60"    # Advanced
    "B6" = "1393"  # TOTAL COURSES
}

foreach ($addr in $levelCounts.Keys) {
    $cell = $wsLevel.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $levelCounts[$addr]
}
